# Swap the contents of column E (codeforiati:category-name) and
# column F (codeforiati:group-code) for every data row in the sheet
# (header row included, so the header labels swap too).
#
# Column F's data values are numeric-looking strings (e.g. "110"), and
# once they move into column E we must keep them as text (matching the
# source file, where every cell -- including these -- is stored as a
# shared string) rather than letting Excel auto-convert them to numbers.
# So format the destination data cells (E2:E235) as Text before writing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$ws.Range("E2:E" + $lastRow).NumberFormat = "@"

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    $eCell.Value2 = $fVal
    $fCell.Value2 = $eVal
}
